# Apply updated cryptocurrency market data values (price & volume columns)
# Also corrects row ordering for a few coins whose relative ranking changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($Worksheet, $CellRef, $Text) {
    $cell = $Worksheet.Range($CellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
}

Set-TextCell $ws "D2" "60.605.66"
Set-TextCell $ws "E2" "  +2.80%  "
Set-TextCell $ws "D3" "2.605.31"
Set-TextCell $ws "E4" "  -2.81%  "
Set-TextCell $ws "D5" "567.21"
Set-TextCell $ws "E5" "  +0.17%  "
Set-TextCell $ws "D6" "142.71"
Set-TextCell $ws "E6" "  +0.54%  "
Set-TextCell $ws "E7" "  -0.33%  "
Set-TextCell $ws "D8" "0.599"
Set-TextCell $ws "E8" "  +1.02%  "
Set-TextCell $ws "D9" "2.630.91"
Set-TextCell $ws "E9" "  +2.06%  "
Set-TextCell $ws "E10" "  -2.09%  "
Set-TextCell $ws "E11" "  +2.96%  "
Set-TextCell $ws "E12" "  -4.66%  "
Set-TextCell $ws "D13" "0.369"
Set-TextCell $ws "E13" "  +7.20%  "
Set-TextCell $ws "D14" "3.070.81"
Set-TextCell $ws "E14" "  +1.53%  "
Set-TextCell $ws "D15" "60.641.78"
Set-TextCell $ws "E15" "  +2.78%  "
Set-TextCell $ws "D16" "23.52"
Set-TextCell $ws "E16" "  +5.22%  "
Set-TextCell $ws "E17" "  +3.46%  "
Set-TextCell $ws "D18" "2.620.31"
Set-TextCell $ws "E18" "  +1.75%  "
Set-TextCell $ws "B19" "Chainlink"
Set-TextCell $ws "C19" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell $ws "D19" "11.26"
Set-TextCell $ws "E19" "  +9.71%  "
Set-TextCell $ws "B20" "Polkadot"
Set-TextCell $ws "C20" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell $ws "D20" "4.65"
Set-TextCell $ws "E20" "  +3.11%  "
Set-TextCell $ws "D21" "347.79"
Set-TextCell $ws "E21" "  +2.88%  "
Set-TextCell $ws "D22" "7.05"
Set-TextCell $ws "E22" "  +12.57%  "
Set-TextCell $ws "D23" "0.999"
Set-TextCell $ws "E23" "  +0.00%  "
Set-TextCell $ws "D24" "0.526"
Set-TextCell $ws "E24" "  +13.78%  "
Set-TextCell $ws "D25" "63.82"
Set-TextCell $ws "E25" "  -0.66%  "
Set-TextCell $ws "E26" "  -0.54%  "
Set-TextCell $ws "E27" "  -1.02%  "
Set-TextCell $ws "D28" "7.69"
Set-TextCell $ws "E28" "  +6.00%  "
Set-TextCell $ws "D29" "0.0₃0789"
Set-TextCell $ws "E29" "  +2.58%  "
Set-TextCell $ws "D30" "1.82"
Set-TextCell $ws "E30" "  +8.65%  "
Set-TextCell $ws "E31" "  -0.14%  "
Set-TextCell $ws "D32" "6.33"
Set-TextCell $ws "E32" "  +4.96%  "
Set-TextCell $ws "D33" "160.47"
Set-TextCell $ws "E33" "  -0.30%  "
Set-TextCell $ws "D34" "19.52"
Set-TextCell $ws "E34" "  +3.40%  "
Set-TextCell $ws "D35" "4.23"
Set-TextCell $ws "E35" "  +6.16%  "
Set-TextCell $ws "D36" "0.968"
Set-TextCell $ws "E36" "  +11.55%  "
Set-TextCell $ws "D37" "1.22"
Set-TextCell $ws "E37" "  +5.56%  "
Set-TextCell $ws "D38" "1.60"
Set-TextCell $ws "E38" "  +8.68%  "
Set-TextCell $ws "D39" "37.75"
Set-TextCell $ws "E39" "  +0.74%  "
Set-TextCell $ws "B40" "Filecoin"
Set-TextCell $ws "C40" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell $ws "D40" "3.80"
Set-TextCell $ws "E40" "  +4.23%  "
Set-TextCell $ws "B41" "SuiNetwork"
Set-TextCell $ws "C41" "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextCell $ws "D41" "0.849"
Set-TextCell $ws "E41" "  -2.00%  "
Set-TextCell $ws "D42" "295.55"
Set-TextCell $ws "E42" "  +0.29%  "
Set-TextCell $ws "D43" "140.30"
Set-TextCell $ws "E43" "  +7.33%  "
Set-TextCell $ws "D44" "0.996"
Set-TextCell $ws "E44" "  -0.36%  "
Set-TextCell $ws "B45" "Mantle"
Set-TextCell $ws "C45" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell $ws "D45" "0.607"
Set-TextCell $ws "E45" "  +2.79%  "
Set-TextCell $ws "B46" "Stellar"
Set-TextCell $ws "C46" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell $ws "D46" "0.0980"
Set-TextCell $ws "E46" "  +0.88%  "
Set-TextCell $ws "D47" "0.0548"
Set-TextCell $ws "E47" "  +2.83%  "
Set-TextCell $ws "B48" "EnergySwap"
Set-TextCell $ws "C48" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell $ws "D48" "19.53"
Set-TextCell $ws "E48" "  +2.38%  "
Set-TextCell $ws "B49" "VeChain"
Set-TextCell $ws "C49" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell $ws "D49" "0.0240"
Set-TextCell $ws "E49" "  +3.54%  "
Set-TextCell $ws "D50" "10.69"
Set-TextCell $ws "E50" "  +0.50%  "
Set-TextCell $ws "B51" "RenderToken"
Set-TextCell $ws "C51" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell $ws "D51" "4.88"
Set-TextCell $ws "E51" "  +8.63%  "
